$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "vk"
$ws.Range("E1").Value = "mail"

$ws.Range("A2").Value = "qwer12"
$ws.Range("B2").Value = "asdc1"
$ws.Range("C2").Value = "tyhe56"
$ws.Range("D2").Value = "hello97"
$ws.Range("E2").Value = "pochta"
